{"js": "// Insert a new empty paragraph right after the \"ng serve \u2013open\" paragraph\n// (i.e. before the first of the already-existing trailing empty paragraphs).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"ng serve \u2013open\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph with text \"ng serve \u2013open\"');\n}\n\ntarget.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new empty paragraph right after the \"ng serve \u2013open\" paragraph\n# (i.e. before the first of the already-existing trailing empty paragraphs).\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"ng serve \u2013open\")\n\nif (-not $found) {\n    throw 'Could not find paragraph with text \"ng serve \u2013open\"'\n}\n\n$range.InsertParagraphAfter()\n"}
